# The workbook's "setlist" table (Table1) had its 4th column header
# renamed from "Event" to "Topic". Updating the header cell's value
# renames the table column in lock-step (ListObjects keep their column
# names driven by the header row cell contents), and also updates the
# shared-string table / D1 cell reference accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Topic"

# Leave the active selection where the author last left it when they
# saved (on the newly renamed header cell's column, one row down).
$ws.Range("D2").Select()
